$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.821.75"
$ws.Range("E2").Value = "  +1.79%  "
$ws.Range("D3").Value = "1.887.53"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.29"
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4725"
$ws.Range("E7").Value = "  +2.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3936"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.67"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08090"
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.028"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.18"
$ws.Range("E12").Value = "  +3.43%  "
$ws.Range("D13").Value = "1.885.79"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.992"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.142"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06735"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001051"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "87.37"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.38"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.008"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "27.834.86"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.532"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.02"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.329"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "2.110.60"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.27"
$ws.Range("E27").Value = "  +3.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.21"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.110"
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.584"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.18"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9839"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09501"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.453"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.624"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.363"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06164"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02267"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.223"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.098"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6014"
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1899"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.34"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.261"
$ws.Range("E44").Value = "  -1.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5720"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.23"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.950"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06921"
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.18"
$ws.Range("E50").Value = "  +4.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000305"
$ws.Range("E51").Value = "  +8.78%  "
